$wb = $excel.ActiveWorkbook

# --- "NamedValue" sheet -> "Concept" sheet ---
$concept = $wb.Worksheets.Item("NamedValue")
$concept.Name = "Concept"

# Rename the synthetic "named_value" URNs / labels to "concept" / "Concept"
$concept.Cells.Replace("named_value", "concept")
$concept.Cells.Replace("Named value", "Concept")

# Drop the now-redundant "title" column (D); "value" (old E) shifts into D
$concept.Columns.Item(4).Delete()

# --- "Image" sheet: its "depicts" column references the old named_value URNs ---
$image = $wb.Worksheets.Item("Image")
$image.Cells.Replace("named_value", "concept")

# --- "Person" sheet: fix one relation URL ---
$person = $wb.Worksheets.Item("Person")
$person.Range("F4").Value = "http://www.wikidata.org/entity/Q7251"
